# Fixed update to excel issue
# - Rename "Requested quantity" header to "Weekly_PO_Qty" on "Weekly Quantity" sheet
# - Rename "Requested quantity" header to "Monthly_PO_Qty" on "Monthly Trend" sheet
# - Add a new "PO Forecast" sheet with ds / PO_Forecast / yhat_lower / yhat_upper data

$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

# --- rename headers -------------------------------------------------------
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- add the new "PO Forecast" sheet at the end ---------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# copy header formatting (bold, centered, bordered) from an existing header
$wsWeekly.Range("A1:B1").Copy() | Out-Null
$wsForecast.Range("A1:D1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$rows = @(
    @(45494.99999999999, 19, -75.56019086677601, 108.2061925380546),
    @(45522.99999999999, 35, -52.97906891223666, 128.2695167095884),
    @(45557.99999999999, 54, -40.73921709358449, 146.7882516468355),
    @(45564.99999999999, 58, -37.36764918883167, 146.0037128375267),
    @(45571.99999999999, 61, -31.46621974030887, 151.4592252052859),
    @(45578.99999999999, 65, -23.03860662847845, 159.9671926038345),
    @(45585.99999999999, 69, -19.63141161062773, 162.7554870318056),
    @(45592.99999999999, 73, -17.45240656778304, 163.9978495723221),
    @(45599.99999999999, 77, -13.92056520371368, 173.540746934993),
    @(45606.99999999999, 81, -12.04043508985082, 178.0749566684471),
    @(45613.99999999999, 84, -2.199377309173539, 172.7392283724523),
    @(45634.99999999999, 96, 5.274994276885536, 187.8544037223732),
    @(45641.99999999999, 100, 12.49960424706407, 192.5033696091974),
    @(45648.99999999999, 104, 16.659159887951, 193.1230053616139),
    @(45655.99999999999, 108, 20.75227856381598, 192.5180995169256),
    @(45662.99999999999, 111, 23.82557792108306, 198.0183054535029),
    @(45669.99999999999, 115, 19.08150236262851, 204.9248404827731),
    @(45676.99999999999, 119, 24.85430951876771, 208.6269143853467),
    @(45683.99999999999, 123, 34.13709110330504, 217.9489603008084),
    @(45690.99999999999, 127, 31.39656850334905, 213.5935542549611)
)

$r = 2
foreach ($row in $rows) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# copy the date-column number format (ds column) from an existing sheet
$wsWeekly.Range("A2").Copy() | Out-Null
$wsForecast.Range("A2:A21").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$wsForecast.Range("A1").Select() | Out-Null

# restore original active sheet/selection (unchanged by this edit)
$wsWeekly.Activate() | Out-Null
$wsWeekly.Range("A1").Select() | Out-Null
